$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.54
$ws.Range("J2").Value = 3
$ws.Range("P2").Value = 1.9
$ws.Range("Q2").Value = 1.87
$ws.Range("G3").Value = 3.1
$ws.Range("I3").Value = 3.05
$ws.Range("J3").Value = 3.45
$ws.Range("Q3").Value = 1.74
$ws.Range("G4").Value = 3.15
$ws.Range("H4").Value = 2.64
$ws.Range("I4").Value = 3.3
$ws.Range("J4").Value = 3.15
$ws.Range("K4").Value = 4.1
$ws.Range("F5").Value = 8.6
$ws.Range("J5").Value = 5.2
$ws.Range("K5").Value = 5.4
$ws.Range("AC5").Value = 12
$ws.Range("AE5").Value = 15
$ws.Range("AF5").Value = 80
$ws.Range("R6").Value = 1.53
$ws.Range("AJ6").Value = 1000
$ws.Range("Q8").Value = 1.75
$ws.Range("R8").Value = 1.5
$ws.Range("S8").Value = 2.9
$ws.Range("T8").Value = 2.26
$ws.Range("AB8").Value = 36
$ws.Range("AF8").Value = 130
$ws.Range("AG8").Value = 48
$ws.Range("AK8").Value = 280
$ws.Range("AL8").Value = 200
$ws.Range("AN8").Value = 450
$ws.Range("N9").Value = 3.85
$ws.Range("P9").Value = 1.99
$ws.Range("Q9").Value = 1.98
$ws.Range("R9").Value = 1.38
$ws.Range("S9").Value = 3.5
$ws.Range("Y9").Value = 12.5
$ws.Range("Z9").Value = 21
$ws.Range("AJ9").Value = 38
$ws.Range("AL9").Value = 40
$ws.Range("AO9").Value = 29
$ws.Range("G10").Value = 4.2
$ws.Range("H10").Value = 1.96
$ws.Range("S10").Value = 2.72
$ws.Range("T10").Value = 1.65
$ws.Range("U10").Value = 2.42
$ws.Range("AG10").Value = 17.5
$ws.Range("AK10").Value = 44
$ws.Range("AL10").Value = 46
$ws.Range("F11").Value = 5
$ws.Range("G11").Value = 5.2
$ws.Range("I11").Value = 1.75
$ws.Range("J11").Value = 4.4
$ws.Range("AE11").Value = 16
$ws.Range("AF11").Value = 44
$ws.Range("AH11").Value = 17
$ws.Range("AL11").Value = 48
$ws.Range("AN11").Value = 40
